$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.564.61'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.98%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.597.77'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -1.75%  '

$ws.Range('E4').Value = '  +0.42%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '208.29'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.30%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.503'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.50%  '

$ws.Range('E7').Value = '  +0.47%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '22.34'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -4.06%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.253'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -1.67%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0593'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.10%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0865'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.58%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.824.88'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.84%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.607.19'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.10%  '

$ws.Range('E14').Value = '  -3.60%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.541'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -3.50%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '63.47'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.74%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '27.562.27'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.99%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '218.73'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -4.63%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.42'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.72%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0697'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -3.20%  '

$ws.Range('E21').Value = '  +0.52%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.22'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.26%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.71'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -3.61%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.99'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.36%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '154.40'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.25%  '

$ws.Range('E26').Value = '  -1.95%  '

$ws.Range('E27').Value = '  +0.45%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.09'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.60%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -3.94%  '

$ws.Range('E31').Value = '  -2.25%  '

$ws.Range('E32').Value = '  -4.30%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.367.70'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.01%  '

$ws.Range('E34').Value = '  -4.32%  '

$ws.Range('E35').Value = '  -2.32%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.972'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -4.69%  '

$ws.Range('E37').Value = '  -0.77%  '

$ws.Range('E38').Value = '  -2.31%  '

$ws.Range('E39').Value = '  -2.55%  '

$ws.Range('E40').Value = '  -4.14%  '

$ws.Range('E41').Value = '  +0.52%  '

$ws.Range('E42').Value = '  -2.66%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.37'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.08%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.79'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.87%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '64.13'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.24%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.734.48'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.98%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.11'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.78%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '88.21'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.28%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0₇0987'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.28%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0973'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -3.77%  '

$ws.Range('E51').Value = '  -0.97%  '

